$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.840.33'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.640.31'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'215.78"
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.50%  '
$ws.Range('D9').Value = "'0.0638"
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').Value = "'19.76"
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('D11').Value = "'0.0795"
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('D12').Value = "'4.27"
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').Value = '1.867.31'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = '1.641.57'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').Value = '0.0₃0769'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').Value = "'63.17"
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '25.874.95'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('D21').Value = "'193.15"
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').Value = "'9.99"
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = "'6.37"
$ws.Range('E23').Value = '  +2.47%  '
$ws.Range('E24').Value = '  +4.32%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').Value = "'142.39"
$ws.Range('E26').Value = '  +2.86%  '
$ws.Range('D27').Value = "'0.122"
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('D29').Value = "'15.57"
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = "'0.0495"
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('D37').Value = '1.133.22'
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = "'0.548"
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'2.53"
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +1.21%  '
$ws.Range('D43').Value = "'100.92"
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('D44').Value = "'0.807"
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('D45').Value = '1.776.47'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('E46').Value = '  +3.97%  '
$ws.Range('D47').Value = "'55.51"
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('E48').Value = '  +6.87%  '
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = "'2.33"
$ws.Range('E51').Value = '  +3.04%  '
